$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Ají" (Chilena(o)) which pushes
# the existing rows 19-59 down to 20-60 (insert, not overwrite).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record's data. The
# constant columns (A, B, C, E, F, G, I, R) carry the same values used by
# every other row in this subset (market/category/quality/classification).
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value = "Bíobío"
$ws.Cells.Item(19, 4).Value = 44526
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100112021
$ws.Cells.Item(19, 7).Value = "Ají"
$ws.Cells.Item(19, 8).Value = "Chilena(o)"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 60
$ws.Cells.Item(19, 11).Value = 50000
$ws.Cells.Item(19, 12).Value = 52000
$ws.Cells.Item(19, 13).Value = 51000
$ws.Cells.Item(19, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(19, 16).Value = 2040
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
